$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing sensor-data rows (columns C:H, rows 2-21) down by one row
# so that row N gets what used to be in row N-1 (N from 21 down to 3).
$cols = @("C", "D", "E", "F", "G", "H")
for ($r = 21; $r -ge 3; $r--) {
    foreach ($col in $cols) {
        $srcCell = $ws.Range($col + ($r - 1))
        $dstCell = $ws.Range($col + $r)
        $dstCell.Value = $srcCell.Value2
    }
}

# New sample data inserted as row 2 (timestamp/label stay as-is, only ax..gz are new)
$ws.Range("C2").Value = -0.647717118263246
$ws.Range("D2").Value = 0.6091025024652482
$ws.Range("E2").Value = -1.090710066258908
$ws.Range("F2").Value = -0.2335032373666763
$ws.Range("G2").Value = -0.1345430761575698
$ws.Range("H2").Value = 0.1078177168965339

# The window is fixed-size, so the former last row (22) is dropped entirely.
$ws.Rows.Item(22).Delete()
